# Update "想去人数" (want-to-go count) figures for two events that appear
# on both the "展览" sheet and the consolidated "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览": row 3 -> 熊喵M动漫嘉年华 (1237 -> 1242), row 4 -> 北极光动漫展 (2729 -> 2732)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 1242
$wsExpo.Range("F4").Value = 2732

# Sheet "全部类型": row 5 -> 熊喵M动漫嘉年华 (1237 -> 1242), row 6 -> 北极光动漫展 (2729 -> 2732)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 1242
$wsAll.Range("F6").Value = 2732
